$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Row 2
$ws.Range("B2").Value = 822.81
$ws.Range("E2").Value = 9177.19
$ws.Range("E2").NumberFormat = "#,##0.00"
$ws.Range("F2").Value = 9177.19
$ws.Range("F2").NumberFormat = "#,##0.00"

# Row 3
$ws.Range("B3").Value = 23.08
$ws.Range("E3").Value = 127.56
$ws.Range("F3").Value = 127.56

# Row 5
$ws.Range("B5").Value = 100
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Make "Summary" the active/selected sheet with a new selection, and move
# selection away from the previously-active "NewLoanInput" sheet.
[void]$ws.Activate()
[void]$ws.Range("K6").Select()
